$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 1786
$ws.Range("L3").Value = 1820
$ws.Range("K4").Value = 1756
$ws.Range("L4").Value = 510
$ws.Range("L5").Value = 110
$ws.Range("L6").Value = 1676
$ws.Range("K7").Value = 27546
$ws.Range("L7").Value = 5902

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 44
$ws.Range("L7").Value = 196
$ws.Range("L8").Value = 363
$ws.Range("L9").Value = 36
$ws.Range("L10").Value = 37
$ws.Range("L11").Value = 107
$ws.Range("L13").Value = 8
$ws.Range("L18").Value = 44
$ws.Range("L19").Value = 169
$ws.Range("L20").Value = 160
$ws.Range("L23").Value = 61
$ws.Range("L25").Value = 30
$ws.Range("L29").Value = 301
$ws.Range("L33").Value = 260
$ws.Range("K36").Value = 357
$ws.Range("L37").Value = 211
$ws.Range("L38").Value = 7
$ws.Range("L47").Value = 38
$ws.Range("L48").Value = 87
$ws.Range("L51").Value = 71
$ws.Range("L52").Value = 126
$ws.Range("L59").Value = 9
$ws.Range("L64").Value = 42
$ws.Range("L65").Value = 114
$ws.Range("L67").Value = 206
$ws.Range("L69").Value = 12
$ws.Range("L72").Value = 24
$ws.Range("L73").Value = 50
$ws.Range("L79").Value = 162
$ws.Range("L85").Value = 305
$ws.Range("L88").Value = 78
$ws.Range("L90").Value = 57
$ws.Range("L96").Value = 54
$ws.Range("L99").Value = 92
$ws.Range("K101").Value = 27546
$ws.Range("L101").Value = 5902

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 11
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 57
$ws.Range("L4").Value = 17
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 196

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 36
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 107

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 91
$ws.Range("L3").Value = 127
$ws.Range("L4").Value = 27
$ws.Range("L5").Value = 6
$ws.Range("L7").Value = 305

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 37
$ws.Range("L7").Value = 126

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 12

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 100
$ws.Range("L3").Value = 122
$ws.Range("L6").Value = 101
$ws.Range("L7").Value = 363

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 90
$ws.Range("L7").Value = 260

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 66
$ws.Range("L7").Value = 211

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L6").Value = 31
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 92

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 60
$ws.Range("L7").Value = 206

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 96
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 83
$ws.Range("L7").Value = 301

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 11
$ws.Range("L7").Value = 87

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 53
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("L3").Value = 3
$ws.Range("L6").Value = 8

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 19
$ws.Range("L3").Value = 24
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 58
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 162

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L2").Value = 15
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 42

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L3").Value = 45
$ws.Range("L7").Value = 160

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L2").Value = 16
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K4").Value = 21
$ws.Range("K7").Value = 357

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 13
$ws.Range("L7").Value = 50

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("L3").Value = 4
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 14
$ws.Range("L6").Value = 15
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 19
$ws.Range("L3").Value = 26
$ws.Range("L7").Value = 78

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 57

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L3").Value = 22
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 24

$ws = $wb.Worksheets.Item("Grant Park")
$ws.Range("L3").Value = 2
$ws.Range("L6").Value = 7
